# Update "想去人数" (want-to-go count) figures across the four sheets
# to the newly scraped values.

$wb = $excel.ActiveWorkbook

# 展览 sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 974
$ws1.Range("F9").Value = 2273
$ws1.Range("F13").Value = 1181
$ws1.Range("F15").Value = 2290
$ws1.Range("F16").Value = 744
$ws1.Range("F17").Value = 16743
$ws1.Range("F20").Value = 598
$ws1.Range("F25").Value = 128
$ws1.Range("F29").Value = 43

# 演出 sheet
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 48

# 本地生活 sheet
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 498

# 全部类型 sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 498
$ws4.Range("F5").Value = 974
$ws4.Range("F9").Value = 2273
$ws4.Range("F15").Value = 1181
$ws4.Range("F20").Value = 2290
$ws4.Range("F21").Value = 744
$ws4.Range("F22").Value = 16743
$ws4.Range("F27").Value = 598
$ws4.Range("F32").Value = 128
$ws4.Range("F40").Value = 43
